$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 143; this shifts the existing rows 143-195
# down to 144-196 and carries their formatting with them.
$ws.Rows.Item(143).Insert()

# Populate the newly inserted row 143 with the new record. Columns that are
# constant across this sheet (A,B,C,E,F,G,H,I,J) as well as K,L,Q,R,T are
# copied from the record that used to occupy row 143 (now row 144); only
# D,M,N,O,P,S carry genuinely new values.
$ws.Cells.Item(143, 1).Value = 7
$ws.Cells.Item(143, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(143, 3).Value = "Ñuble"
$ws.Cells.Item(143, 4).Value = 44809
$ws.Cells.Item(143, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(143, 5).Value = 16
$ws.Cells.Item(143, 6).Value = "Fruta"
$ws.Cells.Item(143, 7).Value = 100102
$ws.Cells.Item(143, 8).Value = "Cítricos"
$ws.Cells.Item(143, 9).Value = 100102004
$ws.Cells.Item(143, 10).Value = "Mandarina"
$ws.Cells.Item(143, 11).Value = "Murcott"
$ws.Cells.Item(143, 12).Value = "Primera"
$ws.Cells.Item(143, 13).Value = 160
$ws.Cells.Item(143, 14).Value = 7500
$ws.Cells.Item(143, 15).Value = 8000
$ws.Cells.Item(143, 16).Value = 7750
$ws.Cells.Item(143, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(143, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(143, 19).Value = 775
$ws.Cells.Item(143, 20).Value = 10
